# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the a5a6e20e-0fcb-4834-b73d-ac720a9a06ff file (row 5 on
# every sheet) to reflect a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G5").Value = "2017-02-09 14:44:42"

# zh-cn sheet: column H = "Latest Handoff Datetime"
$wsZhCn.Range("H5").Value = "2017-02-09 14:44:25"

# de-de sheet: column H = "Latest Handoff Datetime"
$wsDeDe.Range("H5").Value = "2017-02-09 14:44:42"
